# Generate Report for Handback
#
# The localization-status workbook tracks handoff/handback state for each
# localized file. The 385087e5-dd0f-4ddd-87ae-553b940f6a5c.md file has now
# been handed back (and is in sync with en-US), so update its status and
# timestamps on all three sheets (Overview, zh-cn, de-de) and clear the
# stale "not the latest" error detail.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the 385087e5...md file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status

# --- zh-cn sheet: row 3 is the 385087e5...md file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("K3").Value = "2016-08-24 02:48:12"
$zhcn.Range("P3").Value = ""

# --- de-de sheet: row 3 is the 385087e5...md file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("K3").Value = "2016-08-24 02:48:19"
$dede.Range("P3").Value = ""
